# Refresh market-price-derived columns (H:N) in the Leve-profit tables on each
# job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) with the latest values from
# the scheduled market-data runner. Column layout (per table):
#   H currentAveragePrice     I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ             L LevePriceHQ             M LeveProfitNQ
#   N LeveProfitHQ
# LeveProfit cells (M/N) are blank whenever the corresponding Leve price is 0
# (no valid HQ/NQ market data), so some rows gain/lose an N (or M) cell.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 2
$ws.Range("H2").Value = 196
$ws.Range("I2").Value = 61.666668
$ws.Range("J2").Value = 1002
$ws.Range("K2").Value = 61.666668
$ws.Range("L2").Value = 1002
$ws.Range("M2").Value = 51.333332
$ws.Range("N2").Value = -1228

# ALC row 9
$ws.Range("H9").Value = 1876.5834
$ws.Range("I9").Value = 2051.6
$ws.Range("J9").Value = 1001.5
$ws.Range("K9").Value = 2051.6
$ws.Range("L9").Value = 1001.5
$ws.Range("M9").Value = -1882.6
$ws.Range("N9").Value = -1339.5

# ALC row 29
$ws.Range("H29").Value = 8876.4375
$ws.Range("J29").Value = 9441.333000000001
$ws.Range("L29").Value = 28323.999
$ws.Range("N29").Value = -28885.999

# ALC row 33
$ws.Range("H33").Value = 445.57895
$ws.Range("I33").Value = 477.47058
$ws.Range("J33").Value = 174.5
$ws.Range("K33").Value = 477.47058
$ws.Range("L33").Value = 174.5
$ws.Range("M33").Value = -248.47058
$ws.Range("N33").Value = -632.5

# ALC row 38
$ws.Range("H38").Value = 762.5454999999999
$ws.Range("I38").Value = 138.9
$ws.Range("K38").Value = 416.7
$ws.Range("M38").Value = -44.70000000000005

# ALC row 39
$ws.Range("H39").Value = 183.75
$ws.Range("I39").Value = 78
$ws.Range("K39").Value = 234
$ws.Range("M39").Value = 62

# ALC row 58
$ws.Range("H58").Value = 2242.6
$ws.Range("I58").Value = 571.1667
$ws.Range("K58").Value = 1713.5001
$ws.Range("M58").Value = -1563.5001

# ALC row 87
$ws.Range("H87").Value = 79999
$ws.Range("J87").Value = 79999
$ws.Range("L87").Value = 79999
$ws.Range("N87").Value = -82495

# ALC row 88
$ws.Range("H88").Value = 2412.8333
$ws.Range("I88").Value = 1555.125
$ws.Range("J88").Value = 3099
$ws.Range("K88").Value = 1555.125
$ws.Range("L88").Value = 3099
$ws.Range("M88").Value = -1149.125
$ws.Range("N88").Value = -3911

# ALC row 90
$ws.Range("H90").Value = 79999
$ws.Range("J90").Value = 79999
$ws.Range("L90").Value = 239997
$ws.Range("N90").Value = -252477

# ALC row 91
$ws.Range("H91").Value = 2412.8333
$ws.Range("I91").Value = 1555.125
$ws.Range("J91").Value = 3099
$ws.Range("K91").Value = 1555.125
$ws.Range("L91").Value = 3099
$ws.Range("M91").Value = -151.125
$ws.Range("N91").Value = -5907

# ALC row 107
$ws.Range("H107").Value = 6972.1665
$ws.Range("I107").Value = 5499.25
$ws.Range("K107").Value = 5499.25
$ws.Range("M107").Value = -3579.25

# ALC row 125
$ws.Range("H125").Value = 1744.579
$ws.Range("I125").Value = 2485.6667
$ws.Range("J125").Value = 1077.6
$ws.Range("K125").Value = 22371.0003
$ws.Range("L125").Value = 9698.4
$ws.Range("M125").Value = -19911.0003
$ws.Range("N125").Value = -14618.4

# ALC row 132
$ws.Range("H132").Value = 533026
$ws.Range("J132").Value = 10772.6
$ws.Range("L132").Value = 32317.8
$ws.Range("N132").Value = -37377.8

# ALC row 141
$ws.Range("H141").Value = 2524.1428
$ws.Range("I141").Value = 2278.1667
$ws.Range("J141").Value = 4000
$ws.Range("K141").Value = 6834.500100000001
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = -1654.500100000001
$ws.Range("N141").Value = -22360

$ws = $wb.Worksheets.Item("ARM")
# ARM row 25
$ws.Range("H25").Value = 808
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# ARM row 31
$ws.Range("H31").Value = 2306.2
$ws.Range("I31").Value = 1840.2222
$ws.Range("J31").Value = 6500
$ws.Range("K31").Value = 1840.2222
$ws.Range("L31").Value = 6500
$ws.Range("M31").Value = -1546.2222
$ws.Range("N31").Value = -7088

# ARM row 32
$ws.Range("H32").Value = 15227383
$ws.Range("I32").Value = 15156501
$ws.Range("K32").Value = 15156501
$ws.Range("M32").Value = -15156214

# ARM row 35
$ws.Range("H35").Value = 25249.25
$ws.Range("I35").Value = 22666
$ws.Range("J35").Value = 32999
$ws.Range("K35").Value = 22666
$ws.Range("L35").Value = 32999
$ws.Range("M35").Value = -22260
$ws.Range("N35").Value = -33811

# ARM row 61
$ws.Range("H61").Value = 3980.3157
$ws.Range("I61").Value = 4171.5312
$ws.Range("K61").Value = 4171.5312
$ws.Range("M61").Value = -3959.5312

# ARM row 63
$ws.Range("H63").Value = 3896.9167
$ws.Range("J63").Value = 4726.5
$ws.Range("L63").Value = 4726.5
$ws.Range("N63").Value = -6098.5

# ARM row 66
$ws.Range("H66").Value = 3896.9167
$ws.Range("J66").Value = 4726.5
$ws.Range("L66").Value = 23632.5
$ws.Range("N66").Value = -30496.5

# ARM row 74
$ws.Range("H74").Value = 5456.278
$ws.Range("I74").Value = 5515.2144
$ws.Range("K74").Value = 5515.2144
$ws.Range("M74").Value = -4641.2144

# ARM row 77
$ws.Range("H77").Value = 5456.278
$ws.Range("I77").Value = 5515.2144
$ws.Range("K77").Value = 27576.072
$ws.Range("M77").Value = -23208.072

# ARM row 88
$ws.Range("H88").Value = 6186.5835
$ws.Range("J88").Value = 7376.7144
$ws.Range("L88").Value = 7376.7144
$ws.Range("N88").Value = -8188.7144

# ARM row 91
$ws.Range("H91").Value = 6186.5835
$ws.Range("J91").Value = 7376.7144
$ws.Range("L91").Value = 7376.7144
$ws.Range("N91").Value = -10184.7144

# ARM row 122
$ws.Range("H122").Value = 3474.4565
$ws.Range("I122").Value = 2702.3438
$ws.Range("J122").Value = 5239.2856
$ws.Range("K122").Value = 8107.0314
$ws.Range("L122").Value = 15717.8568
$ws.Range("M122").Value = -5657.0314
$ws.Range("N122").Value = -20617.8568

# ARM row 136
$ws.Range("H136").Value = 3980.3157
$ws.Range("I136").Value = 4171.5312
$ws.Range("K136").Value = 12514.5936
$ws.Range("M136").Value = -9964.5936

$ws = $wb.Worksheets.Item("BSM")
# BSM row 37
$ws.Range("H37").Value = 2284.875
$ws.Range("I37").Value = 1278
$ws.Range("J37").Value = 9333
$ws.Range("K37").Value = 1278
$ws.Range("L37").Value = 9333
$ws.Range("M37").Value = -1141
$ws.Range("N37").Value = -9607

$ws = $wb.Worksheets.Item("CRP")
# CRP row 62
$ws.Range("H62").Value = 7261.769
$ws.Range("J62").Value = 6538.8
$ws.Range("L62").Value = 6538.8
$ws.Range("N62").Value = -7786.8

# CRP row 65
$ws.Range("H65").Value = 7261.769
$ws.Range("J65").Value = 6538.8
$ws.Range("L65").Value = 32694
$ws.Range("N65").Value = -38934

# CRP row 122
$ws.Range("H122").Value = 2927.889
$ws.Range("I122").Value = 2471.125
$ws.Range("J122").Value = 3841.4167
$ws.Range("K122").Value = 7413.375
$ws.Range("L122").Value = 11524.2501
$ws.Range("M122").Value = -4963.375
$ws.Range("N122").Value = -16424.2501

# CRP row 138
$ws.Range("H138").Value = 69090.82000000001
$ws.Range("J138").Value = 62222.11
$ws.Range("L138").Value = 62222.11
$ws.Range("N138").Value = -72502.11

$ws = $wb.Worksheets.Item("CUL")
# CUL row 52
$ws.Range("H52").Value = 4000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 4000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 12000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -12532

$ws = $wb.Worksheets.Item("GSM")
# GSM row 102
$ws.Range("H102").Value = 5403.8887
$ws.Range("I102").Value = 4117.5
$ws.Range("K102").Value = 4117.5
$ws.Range("M102").Value = -2495.5

# GSM row 132
$ws.Range("H132").Value = 5235.64
$ws.Range("I132").Value = 5557.6763
$ws.Range("K132").Value = 16673.0289
$ws.Range("M132").Value = -14143.0289

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 6280.5
$ws.Range("I40").Value = 5377.5
$ws.Range("K40").Value = 5377.5
$ws.Range("M40").Value = -5241.5

$ws = $wb.Worksheets.Item("WVR")
# WVR row 136
$ws.Range("H136").Value = 11374344
$ws.Range("I136").Value = 15634920
$ws.Range("J136").Value = 12810.333
$ws.Range("K136").Value = 46904760
$ws.Range("L136").Value = 38430.999
$ws.Range("M136").Value = -46902210
$ws.Range("N136").Value = -43530.999

